$wb = $excel.ActiveWorkbook

# Sheet 1: "t_size=10, mask_p=0.8"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 0
$ws1.Range("G2").Value = 0.8823529411764706
$ws1.Range("H2").Value = 0.8823529411764706
$ws1.Range("G3").Value = 1
$ws1.Range("H3").Value = 1
$ws1.Range("G4").Value = 0.7169811320754716
$ws1.Range("H4").Value = 0.8837209302325582
$ws1.Range("G5").Value = 0.85
$ws1.Range("H5").Value = 0.6538461538461539

# Sheet 2: "t_size=25, mask_p=0.8"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = 0.8888888888888888
$ws2.Range("H2").Value = 0.8888888888888888
$ws2.Range("G3").Value = 0.8181818181818182
$ws2.Range("H3").Value = 0.9
$ws2.Range("G4").Value = 0.7681159420289855
$ws2.Range("H4").Value = 0.6883116883116883
$ws2.Range("G5").Value = 0.7931034482758621
$ws2.Range("H5").Value = 0.7931034482758621

# Sheet 3: "t_size=50, mask_p=0.8"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("G2").Value = 0.7777777777777778
$ws3.Range("H2").Value = 0.7777777777777778
$ws3.Range("G3").Value = 0.8
$ws3.Range("H3").Value = 0.8
$ws3.Range("G4").Value = 0.8095238095238095
$ws3.Range("H4").Value = 0.8947368421052632
$ws3.Range("G5").Value = 0.8823529411764706
$ws3.Range("H5").Value = 0.7894736842105263

# Sheet 4: "t_size=100, mask_p=0.8"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G2").Value = 0.8121546961325967
$ws4.Range("H2").Value = 0.8963414634146342
$ws4.Range("G3").Value = 0.8963414634146342
$ws4.Range("H3").Value = 0.8121546961325967
$ws4.Range("G4").Value = 0.9105691056910568
$ws4.Range("H5").Value = 0.9017857142857144
